# Weekly update: insert a new price record (row 68) for
# "Agrícola del Norte S.A. de Arica - Locoto" and push the existing
# records down by one row (old row 68 -> 69, ..., old row 120 -> 121).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 68; this shifts rows 68:120 down to
# 69:121 (carrying their formatting/styles with them), exactly like
# Excel's own "Insert Sheet Rows" command.
$ws.Rows("68:68").Insert()

# Populate the freshly inserted row 68 with the new weekly record.
$ws.Cells.Item(68, 1).Value = 1
$ws.Cells.Item(68, 2).Value = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(68, 3).Value = "Arica y Parinacota"
$ws.Cells.Item(68, 4).Value = 44778
$ws.Cells.Item(68, 5).Value = 15
$ws.Cells.Item(68, 6).Value = 100112042
$ws.Cells.Item(68, 7).Value = "Locoto"
$ws.Cells.Item(68, 8).Value = "Sin especificar"
$ws.Cells.Item(68, 9).Value = "Segunda"
$ws.Cells.Item(68, 10).Value = 130
$ws.Cells.Item(68, 11).Value = 19000
$ws.Cells.Item(68, 12).Value = 20000
$ws.Cells.Item(68, 13).Value = 19500
$ws.Cells.Item(68, 14).Value = "$/caja 20 kilos"
$ws.Cells.Item(68, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(68, 16).Value = 975
$ws.Cells.Item(68, 17).Value = 20
$ws.Cells.Item(68, 18).Value = "Hortaliza"
